$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56, pushing existing rows 56-182 down to 57-183.
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new observation.
$ws.Cells.Item(56, 1).Value2 = 10
$ws.Cells.Item(56, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(56, 3).Value2 = "La Araucanía"
$ws.Cells.Item(56, 4).Value2 = 44526
$ws.Cells.Item(56, 5).Value2 = 9
$ws.Cells.Item(56, 6).Value2 = 100112039
$ws.Cells.Item(56, 7).Value2 = "Ciboulette"
$ws.Cells.Item(56, 8).Value2 = "Sin especificar"
$ws.Cells.Item(56, 9).Value2 = "Primera"
$ws.Cells.Item(56, 10).Value2 = 20
$ws.Cells.Item(56, 11).Value2 = 5000
$ws.Cells.Item(56, 12).Value2 = 5000
$ws.Cells.Item(56, 13).Value2 = 5000
$ws.Cells.Item(56, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(56, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(56, 16).Value2 = 1667
$ws.Cells.Item(56, 17).Value2 = 3
$ws.Cells.Item(56, 18).Value2 = "Hortaliza"
